$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 26:31 - new colony-count columns D:I (previously blank on these rows)
# ---------------------------------------------------------------------------

# Copy number formats from the existing analogous rows so the new cells pick
# up the same styles (s="1" for D, s="2" for E:I) without minting new xf/font
# records.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D26:D31").PasteSpecial(-4122) | Out-Null
$ws.Range("D56:D61").PasteSpecial(-4122) | Out-Null

$ws.Range("E2:I2").Copy() | Out-Null
$ws.Range("E26:I31").PasteSpecial(-4122) | Out-Null
$ws.Range("E56:I56").PasteSpecial(-4122) | Out-Null

$ws.Range("E2:F2").Copy() | Out-Null
$ws.Range("E57:F61").PasteSpecial(-4122) | Out-Null

$ws.Range("H2").Copy() | Out-Null
$ws.Range("H57:H61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Dilution factor is identical across rows 26:31 (1.0E-6) and 56:61 (1.0E-4)
$ws.Range("D26:D31").Value = 0.000001
$ws.Range("D56:D61").Value = 0.0001

# Purple / White colony counts
$ws.Range("E26").Value = 317
$ws.Range("F26").Value = 124
$ws.Range("E27").Value = 245
$ws.Range("F27").Value = 111
$ws.Range("E28").Value = 302
$ws.Range("F28").Value = 141
$ws.Range("E29").Value = 178
$ws.Range("F29").Value = 150
$ws.Range("E30").Value = 231
$ws.Range("F30").Value = 103
$ws.Range("E31").Value = 324
$ws.Range("F31").Value = 21

# Purple/White CFU per mL + relative abundance, filled down as one shared
# formula family (same pattern as the rest of the column).
$ws.Range("G26:G31").Formula = "=(20*(1/D26))*E26"
$ws.Range("H26:H31").Formula = "=(20*(1/D26))*F26"
$ws.Range("I26:I31").Formula = "=G26/(G26+H26)"

# ---------------------------------------------------------------------------
# Rows 56:61 - new colony-count columns (row 56 gets D:I, rows 57:61 only
# get D:F + H, matching the source data)
# ---------------------------------------------------------------------------
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 262

$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 64
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 189
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 14
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 17
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 162

$ws.Range("G56").Formula = "=(20*(1/D56))*E56"
$ws.Range("H56:H61").Formula = "=(20*(1/D56))*F56"
$ws.Range("I56").Formula = "=G56/(G56+H56)"

# ---------------------------------------------------------------------------
# Sheet view - scroll back to the top and select the newly entered range
# ---------------------------------------------------------------------------
$ws.Range("E56:E61").Select() | Out-Null
